# Rename the first worksheet (was "mIF_template") to "micsss"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "micsss"

# --- Sheet1 ("micsss") data edits -----------------------------------------

# LEAD ORGANIZATION STUDY ID: "E4412" -> 10021 (numeric)
$ws1.Range("C2").Value = 10021

# ASSAY CATEGORY: "Multiplex Immunohistochemistry" -> "Whole Exome Sequencing (WES)"
$ws1.Range("C4").Value = "Whole Exome Sequencing (WES)"

# CELL SEGMENTATION MODEL: "watershed" -> "proprietary"
$ws1.Range("C10").Value = "proprietary"

# PROTOCOL NAME: "Celebi Lab Melanoma Project" -> "Cellebi Lab Melanoma Project"
$ws1.Range("C12").Value = "Cellebi Lab Melanoma Project"

# Row 16 (FOXP3): PR. ANTIBODY INCUBATION TIME was a time value (2 hr) -> literal text "2hr"
$ws1.Range("N16").Value = "2hr"

# Row 17 (CD3): PR. ANTIBODY DILUTION "RTU" -> 9.722222222222221E-2 (same dilution as row 16)
$ws1.Range("M17").Value = 0.09722222222222221
# PR. ANTIBODY INCUBATION TIME -> literal text "2hr" (same convention as row 16)
$ws1.Range("N17").Value = "2hr"
# S. ANTIBODY DILUTION -> literal text "1:20"
$ws1.Range("T17").Value = "1:20"
# S. ANTIBODY INCUBATION TIME: 30 min -> 90 min
$ws1.Range("U17").Value = 0.0625
# AR INCUBATION TIME: 30 min -> 90 min
$ws1.Range("X17").Value = 0.0625

# --- Sheet2 ("Data") data edits --------------------------------------------

# LEAD ORGANIZATION STUDY ID: "E4412" -> 10021 (numeric)
$ws2.Range("C1").Value = 10021

# --- View / selection state -------------------------------------------------
# micsss sheet becomes the active/selected tab, with C21 selected
$ws1.Activate()
$ws1.Range("C21").Select()

# Data sheet selection moves from D10 to C10
$ws2.Range("C10").Select()

# Re-activate micsss so it is the tab shown/selected on open
$ws1.Activate()
